# Generate Report for handback
#
# For both locale sheets (zh-cn, de-de):
#   - Status (col B) for the two content rows changes from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Target File" (col E) and "Latest Handback File" (col F)
#     columns get populated with the same file + hyperlink as the
#     "Source File Name" (col A) and "Latest Handoff File" (col C) columns
#     respectively, since the target of the handback matches the handoff.
#   - The "Latest Handback DateTime" (col G) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/e5bca15541afc9a06016bc46cf0d68f2c7a4e3d4/e2e/1b21140b-1751-4c3d-a04c-404c3db397c4.md", [Type]::Missing, [Type]::Missing, "1b21140b-1751-4c3d-a04c-404c3db397c4.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fcb0b63bd0eb7e71cbd50016e970e04fcb4054d9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/1b21140b-1751-4c3d-a04c-404c3db397c4.d626d21e19a292e8bc71584e8ff55bea664d031c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "1b21140b-1751-4c3d-a04c-404c3db397c4.d626d21e19a292e8bc71584e8ff55bea664d031c.zh-cn.xlf")
$ws.Range("G2").Value = "2016-01-18 02:15:40"

$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5bca15541afc9a06016bc46cf0d68f2c7a4e3d4/e2e/b2545090-b8fe-4df0-8160-19d983cec574.md", [Type]::Missing, [Type]::Missing, "b2545090-b8fe-4df0-8160-19d983cec574.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fcb0b63bd0eb7e71cbd50016e970e04fcb4054d9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/b2545090-b8fe-4df0-8160-19d983cec574.66e19a7e4ba2f67e39f4ab9b09608ee09e41d21f.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "b2545090-b8fe-4df0-8160-19d983cec574.66e19a7e4ba2f67e39f4ab9b09608ee09e41d21f.zh-cn.xlf")
$ws.Range("G3").Value = "2016-01-18 02:15:40"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/e5bca15541afc9a06016bc46cf0d68f2c7a4e3d4/e2e/1b21140b-1751-4c3d-a04c-404c3db397c4.md", [Type]::Missing, [Type]::Missing, "1b21140b-1751-4c3d-a04c-404c3db397c4.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0ca14504393ed73eb21a6239a3cbea0d37d6731/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/1b21140b-1751-4c3d-a04c-404c3db397c4.d626d21e19a292e8bc71584e8ff55bea664d031c.de-de.xlf", [Type]::Missing, [Type]::Missing, "1b21140b-1751-4c3d-a04c-404c3db397c4.d626d21e19a292e8bc71584e8ff55bea664d031c.de-de.xlf")
$ws.Range("G2").Value = "2016-01-18 02:16:00"

$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5bca15541afc9a06016bc46cf0d68f2c7a4e3d4/e2e/b2545090-b8fe-4df0-8160-19d983cec574.md", [Type]::Missing, [Type]::Missing, "b2545090-b8fe-4df0-8160-19d983cec574.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0ca14504393ed73eb21a6239a3cbea0d37d6731/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/b2545090-b8fe-4df0-8160-19d983cec574.66e19a7e4ba2f67e39f4ab9b09608ee09e41d21f.de-de.xlf", [Type]::Missing, [Type]::Missing, "b2545090-b8fe-4df0-8160-19d983cec574.66e19a7e4ba2f67e39f4ab9b09608ee09e41d21f.de-de.xlf")
$ws.Range("G3").Value = "2016-01-18 02:16:00"
